$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 (period 01-01-2021) with revised figures ---
$ws.Cells.Item(74, 2).Value = -1682
$ws.Cells.Item(74, 3).Value = 3004
$ws.Cells.Item(74, 4).Value = 11
$ws.Cells.Item(74, 5).Value = 1333
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 1332

# --- Append new row 75 for period 01-04-2021 ---
# The "date-like" Serie label must land as plain text (matching the rest of
# column A), not get auto-converted into a real Excel date. Writing it as a
# formula that evaluates to text, then collapsing that formula down to its
# static value via copy/paste-values, sidesteps the literal-entry date
# parser entirely and avoids minting a stray number-format style.
$ws.Cells.Item(75, 1).Formula = '="01-04-2021"'
$ws.Cells.Item(75, 1).Copy()
$ws.Cells.Item(75, 1).PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Cells.Item(75, 2).Value = -1457
$ws.Cells.Item(75, 3).Value = 3387
$ws.Cells.Item(75, 4).Value = -35
$ws.Cells.Item(75, 5).Value = 1895
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 1895
